$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.889.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.625.55'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.86%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.78%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  +2.25%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.626.08'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.89%  '

$ws.Range("E10").Value = '  +14.03%  '

$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("E12").Value = '  +1.86%  '

$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000189'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.096.43'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.844.34'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.623.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '366.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.52%  '

$ws.Range("E22").Value = '  -0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.25%  '

$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.05%  '

$ws.Range("E27").Value = '  +4.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '584.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.60%  '

$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("E31").Value = '  +0.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("E33").Value = '  +1.39%  '

$ws.Range("E34").Value = '  -0.73%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("E36").Value = '  -0.88%  '

$ws.Range("E37").Value = '  +0.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.49'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.32%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.372'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.40%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("E42").Value = '  +4.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.00%  '

$ws.Range("E44").Value = '  -0.55%  '

$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.65%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0293'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.05%  '

$ws.Range("E49").Value = '  +1.11%  '

$ws.Range("E50").Value = '  +0.86%  '

$ws.Range("E51").Value = '  +2.24%  '
